$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before E (shifts E..AD to F..AE)
$ws.Columns("E:E").Insert()

# Insert a new row before row 3 (shifts old row 3 -> row 4)
$ws.Rows("3:3").Insert()

Write-Host "after inserts"

# New header for electoral votes column
$ws.Range("E2").Value = "Electoral Votes"

# Electoral votes for Nevada (row4 now)
$ws.Range("E4").Value = 11

# Fix formula text in D4 (parens added by author, value unaffected)
$ws.Range("D4").Formula = "=ROUNDDOWN((1407754/F1)*H1,0)"

# Fill in the new Idaho row (row 3)
$ws.Range("A3").Value = "Idaho"
$ws.Range("B3").Value = "Paulette Jordan"
$ws.Range("C3").Value = "Jim Risch"
$ws.Range("D3").Formula = "=ROUNDDOWN((875000/F1)*H1,0)"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("H3").Formula = "=ROUNDDOWN(458569/F1*H1,0)"
$ws.Range("I3").Value = "Rudy Soto"
$ws.Range("J3").Value = "Russ Fulcher"
$ws.Range("K3").Value = 0.287
$ws.Range("L3").Value = 0.678
$ws.Range("M3").Value = 2
$ws.Range("N3").Formula = "=ROUNDDOWN(391333/F1*H1,0)"
$ws.Range("O3").Value = "C. Aaron Swisher"
$ws.Range("P3").Value = "Mike Simpson"
$ws.Range("Q3").Value = 0.317
$ws.Range("R3").Value = 0.641

# Apply style (format) from the new cellXfs entry (index 3) to the whole Idaho row A3:AE3
$ws.Range("A3:AD3").Style = $ws.Range("A4").Style

